$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: participant names (row 1 header "Team Name" is untouched) ---
$ws.Range("A2").Value  = "Deziree Niki"
$ws.Range("A3").Value  = "Sue Madelon"
$ws.Range("A4").Value  = "Ridge Allissa"
$ws.Range("A5").Value  = "Thad Soan"
$ws.Range("A6").Value  = "Claudine"
$ws.Range("A7").Value  = "Skylar"
$ws.Range("A8").Value  = "Edmé"
$ws.Range("A9").Value  = "Irvin"
$ws.Range("A10").Value = "Flavien"
$ws.Range("A11").Value = "Miranda"
$ws.Range("A12").Value = "Placide"
$ws.Range("A13").Value = "Dory"
$ws.Range("A14").Value = "Stefani"
$ws.Range("A15").Value = "Becky"
$ws.Range("A16").Value = "Roswell"
$ws.Range("A17").Value = "Skye Annabelle"
$ws.Range("A18").Value = "Craig"
$ws.Range("A19").Value = "Carrol"

# --- Column B: locations (row 1 header "Locations" is untouched; now only
#     4 data rows remain instead of the previous L1-L4) ---
$ws.Range("B2").Value = "Merry Mark"
$ws.Range("B3").Value = "Ultman Park"
$ws.Range("B4").Value = "Ken Ross Park"
$ws.Range("B5").Value = "Mulligan Park"

# Rows 6, 8, 9, 11 and 12 previously held a manually-wrapped two-line value
# (ht=28.8) from the old "Team Five/Seven/Eight/Ten/Eleven" text. Re-fit them
# now that shorter single-line names live there so the stale explicit
# height doesn't linger. Rows appended past the old used range (13-19)
# never had a height recorded, so they don't need this.
$ws.Rows.Item(6).AutoFit()
$ws.Rows.Item(8).AutoFit()
$ws.Rows.Item(9).AutoFit()
$ws.Rows.Item(11).AutoFit()
$ws.Rows.Item(12).AutoFit()

# Move the active selection, matching the refreshed view.
[void]$ws.Range("C3").Select()
